# Modify ROS and ROS-lite figures: rename the "Libnoc / Librouting / Libpower"
# library labels to "mppa_noc / mppa_routing / mppa_power" in the
# "Low-Level Library" flow-chart shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape by name/id rather than a hard-coded index, in case shape
# ordering differs from what we inspected.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Id -eq 16) {
        $shp = $s.Shapes.Item($i)
        break
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(15)
}

$tr = $shp.TextFrame.TextRange

# Original text: "Libnoc, Librouting, Libpower"
#   chars 1-6   "Libnoc"      -> becomes two runs: "m" + "ppa_noc"
#   chars 7-8   ", "          -> unchanged
#   chars 9-18  "Librouting"  -> becomes "mppa_routing"
#   chars 19-20 ", "          -> unchanged
#   chars 21-28 "Libpower"    -> becomes "mppa_power"

# Force a run split right after the leading "L" of "Libnoc" (between char 1
# and chars 2-6) so the new text can be written into two distinct runs.
$tr.Characters(2, 5).Font.Name = "Segoe UI"

# Apply replacements from right to left so earlier character offsets stay valid.
$tr.Characters(21, 8).Text = "mppa_power"
$tr.Characters(9, 10).Text = "mppa_routing"
$tr.Characters(2, 5).Text = "ppa_noc"
$tr.Characters(1, 1).Text = "m"
